$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$ws.Range("D2").Value = 1490.67
$ws.Range("E2").Value = -1490.67

$ws.Range("D4").Value = 2402.35
$ws.Range("E4").Value = 11320.99
$ws.Range("F4").Value = 0.1750557808813306
